# BeerTaste.xlsx: add a row under the big title on the TastersSchema sheet
# holding two labelled boxes "Epost:" and "Fodselsaar:" (matching the
# existing Navn/Produsent/... schema-header row style), and leave the
# workbook positioned on that sheet (as the author was, while editing it).

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Tasters")
$ws3 = $wb.Worksheets.Item("TastersSchema")

# Insert a new blank row above the old row 2 (the schema header row),
# pushing the header row and the data-entry row down by one.
$ws3.Rows.Item(2).Insert()

# --- Row 1 (title row): box the whole row with a thin outline, like the
# header/data rows below it already do per-cell (left/mid/mid/.../right).
$ws3.Range("A1:G1").BorderAround(1, 2, -4105, $null)

# --- Row 2 (new row): two boxed labels, "Epost:" spanning A:D and
# "Fodselsaar:" spanning E:G, in the same 18pt font used for these labels.
$ws3.Range("A2").Value = "Epost:"
$ws3.Range("E2").Value = "Fødselsår:"

$box1 = $ws3.Range("A2:D2")
$box1.BorderAround(1, 2, -4105, $null)
$box1.Font.Size = 18
$box1.Font.Bold = $false

$box2 = $ws3.Range("E2:G2")
$box2.BorderAround(1, 2, -4105, $null)
$box2.Font.Size = 18
$box2.Font.Bold = $false

# --- Restore each sheet's own last selection, then land on TastersSchema
# (the sheet that was active when the author saved) with B4 selected.
$ws2.Activate()
$ws2.Range("B3").Select()

$ws3.Activate()
$ws3.Range("B4").Select()
